# Append 4 new DHFR mutant kinetics rows (17-20) to the "kinetics_dhfr" table/sheet,
# matching the data added in the target revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kinetics_dhfr")
$lo = $ws.ListObjects.Item(1)

# New data, column-by-column so new shared strings get interned in the same
# order as the authored workbook (all four mutation names, then the four
# "This work" reference cells).
$mutations = @("E17V", "I5K", "V13H", "M20Q")
$kcat      = @(0.4, 2.1, 0.7, 2.85)
$km        = @(1.2, 62, 1.9, 3)
$kcatStd   = @(0.1, 0.32, 0.05, 0.68)
$kmStd     = @(0.58, 1.4, 0.1, 0.57)
$reference = @("This work", "This work", "This work", "This work")

$startRow = 17

for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $mutations[$i]
}
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $kcat[$i]
}
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $km[$i]
}
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $kcatStd[$i]
}
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($startRow + $i, 5).Value = $kmStd[$i]
}
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($startRow + $i, 6).Value = $reference[$i]
}

# Match formatting used by the rest of the table: column A (mutation names)
# uses the "General" number-format style, columns B-F use the plain Arial
# style already applied to the existing data rows.
$ws.Range("A16").Copy()
$ws.Range("A17:A20").PasteSpecial(-4122) | Out-Null
$ws.Range("B16:F16").Copy()
$ws.Range("B17:F20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Grow the table/list-object so the new rows become part of it.
$lo.Resize($ws.Range("A1:F20"))

# The query-table backed "ExternalData_1" name tracks the table body; move it
# out to the new last row.
$wb.Names.Item("kinetics_dhfr!ExternalData_1").RefersTo = "=kinetics_dhfr!`$A`$1:`$E`$20"

# Match the final cursor position recorded in the authored workbook.
$ws.Range("F23").Select()
